$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update summary values
$ws.Range("B4").Value = "8"

# Row 6 (child 0)
$ws.Range("B6").Value = "16"
$ws.Range("C6").Value = "Collette  "
$ws.Range("D6").Value = "Billi  "
$ws.Range("E6").Value = "9.04,-9.54"
$ws.Range("F6").Value = "Elias(mother): 0578741979"
$ws.Range("H6").Value = "44.0"

# Row 7 (child 1)
$ws.Range("B7").Value = "4"
$ws.Range("C7").Value = "Francisca  "
$ws.Range("D7").Value = "Stevie  "
$ws.Range("E7").Value = "8.52,-5.98"
$ws.Range("F7").Value = "Bernardine(mother): 0561339273"
$ws.Range("G7").Value = "7:05:00"
$ws.Range("H7").Value = "39.0"

# Row 8 (child 2)
$ws.Range("B8").Value = "5"
$ws.Range("C8").Value = "Patti  "
$ws.Range("D8").Value = "Lavenia  "
$ws.Range("E8").Value = "9.35,-5.18"
$ws.Range("F8").Value = "Jennell(mother): 0503029941"
$ws.Range("G8").Value = "7:07:00"
$ws.Range("H8").Value = "37.0"

# Row 9 (child 3)
$ws.Range("B9").Value = "14"
$ws.Range("C9").Value = "Lorinda  "
$ws.Range("D9").Value = "Tyron  "
$ws.Range("E9").Value = "5.68,-4.32"
$ws.Range("F9").Value = "Teresa(grandmother): 0558587699"
$ws.Range("G9").Value = "7:12:00"
$ws.Range("H9").Value = "32.0"

# Row 10 (child 4)
$ws.Range("B10").Value = "18"
$ws.Range("C10").Value = "Kandis  "
$ws.Range("D10").Value = "Zulma  "
$ws.Range("E10").Value = "8.28,-3.72"
$ws.Range("F10").Value = "Kylie(mother): 0575413269"
$ws.Range("G10").Value = "7:15:00"
$ws.Range("H10").Value = "29.0"

# Row 11 (child 5)
$ws.Range("B11").Value = "3"
$ws.Range("C11").Value = "Alexia  "
$ws.Range("D11").Value = "Ramonita  "
$ws.Range("E11").Value = "9.12,0.07"
$ws.Range("F11").Value = "Han(father): 0567537032"
$ws.Range("G11").Value = "7:20:00"
$ws.Range("H11").Value = "24.0"

# Row 12 (child 6)
$ws.Range("B12").Value = "6"
$ws.Range("C12").Value = "Ema  "
$ws.Range("D12").Value = "Ardell  "
$ws.Range("E12").Value = "8.06,7.39"
$ws.Range("F12").Value = "Carley(grandmother): 0533587167"
$ws.Range("G12").Value = "7:30:00"
$ws.Range("H12").Value = "14.0"

# Row 13 (child 7)
$ws.Range("B13").Value = "20"
$ws.Range("C13").Value = "Ron"
$ws.Range("D13").Value = "Cohen"
$ws.Range("E13").Value = "6.33,5.28"
$ws.Range("F13").Value = "Bernardine(mother): 0576270618"
$ws.Range("G13").Value = "7:34:00"
$ws.Range("H13").Value = "10.0"

# Row 14 becomes the "school" row (previously row 15); H14 cleared
$ws.Range("A14").Value = "school"
$ws.Range("B14").Value = "3"
$ws.Range("C14").Value = "Ironiah"
$ws.Range("D14").Value = "mySchool"
$ws.Range("E14").Value = "0,0"
$ws.Range("F14").Value = "Shir(secretary): 0523345098"
$ws.Range("G14").Value = "7:44:00"
$ws.Range("H14").ClearContents()

# Row 15 becomes the "cost" row (previously row 16); C15:G15 cleared
$ws.Range("A15").Value = "cost"
$ws.Range("B15").Value = "25"
$ws.Range("C15:G15").ClearContents()

# Row 16 becomes the "time" row (previously row 17); C16:G16 cleared
$ws.Range("A16").Value = "time"
$ws.Range("B16").Value = "44.0"
$ws.Range("C16:G16").ClearContents()

# Delete the now-obsolete former row 17 (entire row, since data shifted up by one)
$ws.Rows.Item(17).Delete()
